# Updated symbol list (coinranking.com crypto snapshot refresh).
#
# The source feed re-ran and produced new Price / Volume(1h) figures for
# most rows, plus a one-row cyclic re-rank among rows 6-18 (each coin's
# name+link moved up a row, with GateToken wrapping from row 6 to row 18)
# before the refreshed Price/Volume values were written in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin / Link columns that re-ranked (rows 6-18) -----------------------
$textUpdates = @{
    "B6" = "FTXToken"
    "C6" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "B7" = "MXToken"
    "C7" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "B8" = "BTSEToken"
    "C8" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "B9" = "LiechtensteinCryptoassetsExchange"
    "C9" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "B10" = "WazirX"
    "C10" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "B11" = "MandalaExchangeToken"
    "C11" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "B12" = "BitrueCoin"
    "C12" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "B13" = "BitMartToken"
    "C13" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "B14" = "BitForexToken"
    "C14" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "B15" = "TigerCash"
    "C15" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "B16" = "UpBots"
    "C16" = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
    "B17" = "LEO"
    "C17" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "B18" = "GateToken"
    "C18" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
}

foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# --- Price / Volume(1h) columns --------------------------------------------
# These new values are plain numeric- or percentage-looking strings
# ("320.29", "3.70%", ...). The source sheet stores them as literal text,
# so force each target cell to Text format before assigning, then restore
# the default "General" format -- this stops Excel's normal on-entry type
# inference from silently turning them into real numbers/percentages.
$numericUpdates = @{
    "D2" = "320.29"
    "E2" = "3.70%"
    "D3" = "41.37"
    "E3" = "1.29%"
    "E4" = "2.70%"
    "D5" = "0.07749"
    "E5" = "1.70%"
    "D6" = "1.770"
    "E6" = "10.24%"
    "D7" = "0.9439"
    "E7" = "3.97%"
    "D8" = "2.425"
    "E8" = "-1.58%"
    "D9" = "0.1244"
    "E9" = "-2.76%"
    "D10" = "0.1872"
    "E10" = "3.73%"
    "D11" = "0.09196"
    "E11" = "1.77%"
    "D12" = "0.04310"
    "E12" = "0.24%"
    "D13" = "0.1050"
    "E13" = "0.66%"
    "D14" = "0.001291"
    "E14" = "2.72%"
    "D15" = "0.006031"
    "E15" = "6.26%"
    "D16" = "0.007491"
    "E16" = "1,897.31%"
    "D17" = "3.342"
    "E17" = "-0.09%"
    "D18" = "4.343"
    "E18" = "1.39%"
    "D19" = "0.3334"
    "E19" = "0.60%"
    "D20" = "7.759"
    "E20" = "12.15%"
    "D21" = "0.1355"
    "E21" = "-2.76%"
    "D23" = "0.04038"
    "E23" = "-0.14%"
    "E24" = "-0.09%"
    "E25" = "1.79%"
    "D26" = "0.0001271"
    "E26" = "-0.09%"
    "D38" = "0.02547"
    "E38" = "5.47%"
    "D39" = "0.05345"
    "E39" = "2.46%"
    "D40" = "0.007778"
    "E40" = "-0.80%"
    "D41" = "0.1317"
    "E41" = "1.38%"
    "D42" = "0.007068"
    "E42" = "3.83%"
    "D43" = "0.001994"
    "E43" = "3.14%"
    "D44" = "0.008280"
    "E44" = "12.83%"
    "D45" = "0.3176"
    "E45" = "-5.25%"
    "D46" = "0.00006690"
    "E46" = "-3.02%"
    "E47" = "0.03%"
    "D48" = "0.2012"
    "E48" = "53.00%"
    "D49" = "0.004208"
    "E49" = "40.15%"
    "E50" = "0.03%"
    "D51" = "0.0002004"
    "E51" = "0.03%"
}

foreach ($addr in $numericUpdates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $numericUpdates.Keys) {
    $ws.Range($addr).Value = $numericUpdates[$addr]
}
foreach ($addr in $numericUpdates.Keys) {
    $ws.Range($addr).NumberFormat = "General"
}
